$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("IAM&CSV Standard price list")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 29
$r = $excel.ActiveWindow.ScrollRow
Write-Host "ScrollRow after set = $r"
